$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 465, shifting the existing rows 465-484 down to 466-485.
$ws.Rows.Item(465).Insert()

# Fill the freshly inserted row 465 with new weekly data (same Mercado/Region/Categoria
# metadata as its neighbours, new Fecha/Volumen/Precio/Origen figures).
$ws.Cells.Item(465, 1).Value = 10
$ws.Cells.Item(465, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(465, 3).Value = "La Araucanía"
$ws.Cells.Item(465, 4).Value = 45075
$ws.Cells.Item(465, 5).Value = 9
$ws.Cells.Item(465, 6).Value = 100112044
$ws.Cells.Item(465, 7).Value = "Perejil"
$ws.Cells.Item(465, 8).Value = "Sin especificar"
$ws.Cells.Item(465, 9).Value = "Primera"
$ws.Cells.Item(465, 10).Value = 45
$ws.Cells.Item(465, 11).Value = 3700
$ws.Cells.Item(465, 12).Value = 3700
$ws.Cells.Item(465, 13).Value = 3700
$ws.Cells.Item(465, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(465, 15).Value = "Región Metropolitana"
$ws.Cells.Item(465, 16).Value = 1233
$ws.Cells.Item(465, 17).Value = 3
$ws.Cells.Item(465, 18).Value = "Hortaliza"

# Append a brand-new final row 485, carrying forward the data that used to live
# in row 484 before the insert shifted everything down by one.
$ws.Cells.Item(485, 1).Value = 10
$ws.Cells.Item(485, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(485, 3).Value = "La Araucanía"
$ws.Cells.Item(485, 4).Value = 44910
$ws.Cells.Item(485, 5).Value = 9
$ws.Cells.Item(485, 6).Value = 100112044
$ws.Cells.Item(485, 7).Value = "Perejil"
$ws.Cells.Item(485, 8).Value = "Sin especificar"
$ws.Cells.Item(485, 9).Value = "Primera"
$ws.Cells.Item(485, 10).Value = 50
$ws.Cells.Item(485, 11).Value = 4000
$ws.Cells.Item(485, 12).Value = 4000
$ws.Cells.Item(485, 13).Value = 4000
$ws.Cells.Item(485, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(485, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(485, 16).Value = 1333
$ws.Cells.Item(485, 17).Value = 3
$ws.Cells.Item(485, 18).Value = "Hortaliza"

# Match the number format used on the other date cells in column D.
$ws.Cells.Item(485, 4).NumberFormat = $ws.Cells.Item(484, 4).NumberFormat
